# Apply the "Trip Type Table" edit:
#  - Add a new row (Code=0, Type="Unknown") to the data range / table
#  - Resize Table1 to cover the new row
#  - Resize column B to fit the new, wider "Unknown" label (remove bestFit autosize)
#  - Move the active selection to G17 (matches the author's saved selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row -------------------------------------------------
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Unknown"

# --- Grow the Excel Table (ListObject) to include the new row -------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B4"))

# --- Adjust column B width (no longer auto "best fit") --------------------
$ws.Columns.Item(2).ColumnWidth = 15.9166666666667

# --- Update the saved selection to G17, like in the authored workbook -----
$ws.Range("G17").Select() | Out-Null
